$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" between "2021-Q1" and "总计"
# ---------------------------------------------------------------------
$wsQ1_2021 = $wb.Worksheets.Item("2021-Q1")
$wsNew = $wb.Worksheets.Add($null, $wsQ1_2021)
$wsNew.Name = "2022-Q1"

# NOTE: worksheet references are positional, not stable identities - any
# sheet captured before the Add() above may now resolve to the wrong
# (shifted) sheet, so re-resolve "总计" by name only after the insert.
$wsTotal = $wb.Worksheets.Item("总计")

# Copy the header-row / index-column formatting from the "2021-Q1" sheet so
# the new sheet matches the look of the other per-quarter sheets.
$wsQ1_2021.Range("B1:H1").Copy()
$wsNew.Range("B1:H1").PasteSpecial(-4122)
$wsQ1_2021.Range("A2:A3").Copy()
$wsNew.Range("A2:A3").PasteSpecial(-4122)

# Headers
$wsNew.Range("B1").Value = "基金代码"
$wsNew.Range("C1").Value = "基金名称"
$wsNew.Range("D1").Value = "基金规模"
$wsNew.Range("E1").Value = "股票总仓位"
$wsNew.Range("F1").Value = "仓位占比"
$wsNew.Range("G1").Value = "持有市值(亿元)"
$wsNew.Range("H1").Value = "仓位排名"

# Row 2
$wsNew.Range("A2").Value = 0
$wsNew.Range("B2").Value = "'013776"
$wsNew.Range("C2").Value = "中泰兴为价值精选混合A"
$wsNew.Range("D2").Value = "'20.31"
$wsNew.Range("E2").Value = "'85.34"
$wsNew.Range("F2").Value = "'4.30"
$wsNew.Range("G2").Value = "'0.8733"
$wsNew.Range("H2").Value = 6

# Row 3
$wsNew.Range("A3").Value = 1
$wsNew.Range("B3").Value = "'013777"
$wsNew.Range("C3").Value = "中泰兴为价值精选混合C"
$wsNew.Range("D3").Value = "'8.71"
$wsNew.Range("E3").Value = "'85.34"
$wsNew.Range("F3").Value = "'4.30"
$wsNew.Range("G3").Value = "'0.3745"
$wsNew.Range("H3").Value = 6

# ---------------------------------------------------------------------
# 2. Add a new top data row to "总计" for the 2022-Q1 summary, pushing
#    the existing rows down.
# ---------------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert()

$wsTotal.Range("A3:A4").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

$wsTotal.Range("B2:D2").ClearFormats()

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 1.25

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2
